$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.379.74"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.641.76"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.85"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.55"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "2.640.57"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("E10").Value = "  +7.90%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.19"
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "3.122.50"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "68.197.18"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "2.636.39"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.48"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.51"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.62"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.85"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("D29").Value = "2.781.18"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "573.74"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.18"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.43"
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.60"
$ws.Range("E37").Value = "  +4.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.86"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.41"
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.374"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.90"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.43"
$ws.Range("D43").Value = "0.0₆0337"
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.73"
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.67"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.67"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.97"
$ws.Range("E51").Value = "  -0.48%  "

Write-Host "Applied cryptos update"
